$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'WI'
$ws.Cells.Item(2, 2).Value = '$1.00 Games'
$ws.Cells.Item(2, 3).Value = '2 For The Money'
$ws.Cells.Item(2, 4).Value = 2078
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = '''2019-03-12'
$ws.Cells.Item(3, 1).Value = 'WI'
$ws.Cells.Item(3, 2).Value = '$1.00 Games'
$ws.Cells.Item(3, 3).Value = 'Food Series'
$ws.Cells.Item(3, 4).Value = 629
$ws.Cells.Item(3, 5).Value = 17
$ws.Cells.Item(3, 6).Value = '''2019-03-12'
$ws.Cells.Item(4, 1).Value = 'WI'
$ws.Cells.Item(4, 2).Value = '$1.00 Games'
$ws.Cells.Item(4, 3).Value = 'Go Scratch Go!'
$ws.Cells.Item(4, 4).Value = 2103
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = '''2019-03-12'
$ws.Cells.Item(5, 1).Value = 'WI'
$ws.Cells.Item(5, 2).Value = '$1.00 Games'
$ws.Cells.Item(5, 3).Value = 'A Latte Cash '
$ws.Cells.Item(5, 4).Value = 2140
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = '''2019-03-12'
$ws.Cells.Item(6, 1).Value = 'WI'
$ws.Cells.Item(6, 2).Value = '$1.00 Games'
$ws.Cells.Item(6, 3).Value = 'Winter Time Winnings'
$ws.Cells.Item(6, 4).Value = 2117
$ws.Cells.Item(6, 5).Value = 83
$ws.Cells.Item(6, 6).Value = '''2019-03-12'
$ws.Cells.Item(7, 1).Value = 'WI'
$ws.Cells.Item(7, 2).Value = '$1.00 Games'
$ws.Cells.Item(7, 3).Value = '5X The Money'
$ws.Cells.Item(7, 4).Value = 2148
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = '''2019-03-12'
$ws.Cells.Item(8, 1).Value = 'WI'
$ws.Cells.Item(8, 2).Value = '$1.00 Games'
$ws.Cells.Item(8, 3).Value = 'Cheese Series'
$ws.Cells.Item(8, 4).Value = 579
$ws.Cells.Item(8, 5).Value = 24
$ws.Cells.Item(8, 6).Value = '''2019-03-12'
$ws.Cells.Item(9, 1).Value = 'WI'
$ws.Cells.Item(9, 2).Value = '$1.00 Games'
$ws.Cells.Item(9, 3).Value = 'Blackjack'
$ws.Cells.Item(9, 4).Value = 2132
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = '''2019-03-12'
$ws.Cells.Item(10, 1).Value = 'WI'
$ws.Cells.Item(10, 2).Value = '$1.00 Games'
$ws.Cells.Item(10, 3).Value = 'Fast $50''s'
$ws.Cells.Item(10, 4).Value = 2072
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = '''2019-03-12'
$ws.Cells.Item(11, 1).Value = 'WI'
$ws.Cells.Item(11, 2).Value = '$1.00 Games'
$ws.Cells.Item(11, 3).Value = '3-2-Won!'
$ws.Cells.Item(11, 4).Value = 2073
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = '''2019-03-12'
$ws.Cells.Item(12, 1).Value = 'WI'
$ws.Cells.Item(12, 2).Value = '$1.00 Games'
$ws.Cells.Item(12, 3).Value = 'Baseball'
$ws.Cells.Item(12, 4).Value = 2077
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = '''2019-03-12'
$ws.Cells.Item(13, 1).Value = 'WI'
$ws.Cells.Item(13, 2).Value = '$1.00 Games'
$ws.Cells.Item(13, 3).Value = 'red white blue'
$ws.Cells.Item(13, 4).Value = 2076
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = '''2019-03-12'
$ws.Cells.Item(14, 1).Value = 'WI'
$ws.Cells.Item(14, 2).Value = '$1.00 Games'
$ws.Cells.Item(14, 3).Value = 'Blackjack Extra'
$ws.Cells.Item(14, 4).Value = 2068
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = '''2019-03-12'
$ws.Cells.Item(15, 1).Value = 'WI'
$ws.Cells.Item(15, 2).Value = '$1.00 Games'
$ws.Cells.Item(15, 3).Value = 'Blackjack Tripler'
$ws.Cells.Item(15, 4).Value = 2108
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = '''2019-03-12'
$ws.Cells.Item(16, 1).Value = 'WI'
$ws.Cells.Item(16, 2).Value = '$2.00 Games'
$ws.Cells.Item(16, 3).Value = 'Lucky Pairs'
$ws.Cells.Item(16, 4).Value = 2062
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = '''2019-03-12'
$ws.Cells.Item(17, 1).Value = 'WI'
$ws.Cells.Item(17, 2).Value = '$2.00 Games'
$ws.Cells.Item(17, 3).Value = '10X The Money'
$ws.Cells.Item(17, 4).Value = 2149
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = '''2019-03-12'
$ws.Cells.Item(18, 1).Value = 'WI'
$ws.Cells.Item(18, 2).Value = '$2.00 Games'
$ws.Cells.Item(18, 3).Value = 'Triple Win'
$ws.Cells.Item(18, 4).Value = 2024
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = '''2019-03-12'
$ws.Cells.Item(19, 1).Value = 'WI'
$ws.Cells.Item(19, 2).Value = '$2.00 Games'
$ws.Cells.Item(19, 3).Value = 'Easy as 1-2-3'
$ws.Cells.Item(19, 4).Value = 2144
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = '''2019-03-12'
$ws.Cells.Item(20, 1).Value = 'WI'
$ws.Cells.Item(20, 2).Value = '$2.00 Games'
$ws.Cells.Item(20, 3).Value = 'Fish & Chips '
$ws.Cells.Item(20, 4).Value = 2096
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = '''2019-03-12'
$ws.Cells.Item(21, 1).Value = 'WI'
$ws.Cells.Item(21, 2).Value = '$2.00 Games'
$ws.Cells.Item(21, 3).Value = 'Crossword'
$ws.Cells.Item(21, 4).Value = 2050
$ws.Cells.Item(21, 5).Value = 12
$ws.Cells.Item(21, 6).Value = '''2019-03-12'
$ws.Cells.Item(22, 1).Value = 'WI'
$ws.Cells.Item(22, 2).Value = '$2.00 Games'
$ws.Cells.Item(22, 3).Value = 'Luck'
$ws.Cells.Item(22, 4).Value = 2099
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 6).Value = '''2019-03-12'
$ws.Cells.Item(23, 1).Value = 'WI'
$ws.Cells.Item(23, 2).Value = '$2.00 Games'
$ws.Cells.Item(23, 3).Value = 'Farm Fresh Cash'
$ws.Cells.Item(23, 4).Value = 2031
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = '''2019-03-12'
$ws.Cells.Item(24, 1).Value = 'WI'
$ws.Cells.Item(24, 2).Value = '$2.00 Games'
$ws.Cells.Item(24, 3).Value = '10 Times Lucky '
$ws.Cells.Item(24, 4).Value = 2129
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = '''2019-03-12'
$ws.Cells.Item(25, 1).Value = 'WI'
$ws.Cells.Item(25, 2).Value = '$2.00 Games'
$ws.Cells.Item(25, 3).Value = 'Cash In A Flash'
$ws.Cells.Item(25, 4).Value = 2101
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = '''2019-03-12'
$ws.Cells.Item(26, 1).Value = 'WI'
$ws.Cells.Item(26, 2).Value = '$2.00 Games'
$ws.Cells.Item(26, 3).Value = 'Tyrannosaurus BUCKS'
$ws.Cells.Item(26, 4).Value = 2102
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = '''2019-03-12'
$ws.Cells.Item(27, 1).Value = 'WI'
$ws.Cells.Item(27, 2).Value = '$2.00 Games'
$ws.Cells.Item(27, 3).Value = 'Money Quest'
$ws.Cells.Item(27, 4).Value = 2030
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = '''2019-03-12'
$ws.Cells.Item(28, 1).Value = 'WI'
$ws.Cells.Item(28, 2).Value = '$2.00 Games'
$ws.Cells.Item(28, 3).Value = 'Cash X10'
$ws.Cells.Item(28, 4).Value = 2085
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = '''2019-03-12'
$ws.Cells.Item(29, 1).Value = 'WI'
$ws.Cells.Item(29, 2).Value = '$2.00 Games'
$ws.Cells.Item(29, 3).Value = 'Joker''s Wild'
$ws.Cells.Item(29, 4).Value = 2063
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = '''2019-03-12'
$ws.Cells.Item(30, 1).Value = 'WI'
$ws.Cells.Item(30, 2).Value = '$2.00 Games'
$ws.Cells.Item(30, 3).Value = '10X the Money'
$ws.Cells.Item(30, 4).Value = 2060
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = '''2019-03-12'
$ws.Cells.Item(31, 1).Value = 'WI'
$ws.Cells.Item(31, 2).Value = '$2.00 Games'
$ws.Cells.Item(31, 3).Value = 'Kitty Cash Doggy Dough'
$ws.Cells.Item(31, 4).Value = 2095
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = '''2019-03-12'
$ws.Cells.Item(32, 1).Value = 'WI'
$ws.Cells.Item(32, 2).Value = '$2.00 Games'
$ws.Cells.Item(32, 3).Value = 'Love To Win'
$ws.Cells.Item(32, 4).Value = 2059
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = '''2019-03-12'
$ws.Cells.Item(33, 1).Value = 'WI'
$ws.Cells.Item(33, 2).Value = '$2.00 Games'
$ws.Cells.Item(33, 3).Value = '100 Bucks'
$ws.Cells.Item(33, 4).Value = 2029
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = '''2019-03-12'
$ws.Cells.Item(34, 1).Value = 'WI'
$ws.Cells.Item(34, 2).Value = '$3.00 Games'
$ws.Cells.Item(34, 3).Value = 'Toad-ally Awesome Crossword'
$ws.Cells.Item(34, 4).Value = 2163
$ws.Cells.Item(34, 5).Value = 3
$ws.Cells.Item(34, 6).Value = '''2019-03-12'
$ws.Cells.Item(35, 1).Value = 'WI'
$ws.Cells.Item(35, 2).Value = '$3.00 Games'
$ws.Cells.Item(35, 3).Value = 'Block Party Bingo'
$ws.Cells.Item(35, 4).Value = 2120
$ws.Cells.Item(35, 5).Value = 1
$ws.Cells.Item(35, 6).Value = '''2019-03-12'
$ws.Cells.Item(36, 1).Value = 'WI'
$ws.Cells.Item(36, 2).Value = '$3.00 Games'
$ws.Cells.Item(36, 3).Value = 'Twisty Treasures'
$ws.Cells.Item(36, 4).Value = 2152
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = '''2019-03-12'
$ws.Cells.Item(37, 1).Value = 'WI'
$ws.Cells.Item(37, 2).Value = '$3.00 Games'
$ws.Cells.Item(37, 3).Value = 'Cashingo'
$ws.Cells.Item(37, 4).Value = 2053
$ws.Cells.Item(37, 5).Value = 1
$ws.Cells.Item(37, 6).Value = '''2019-03-12'
$ws.Cells.Item(38, 1).Value = 'WI'
$ws.Cells.Item(38, 2).Value = '$3.00 Games'
$ws.Cells.Item(38, 3).Value = 'S''More Slingo'
$ws.Cells.Item(38, 4).Value = 2127
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = '''2019-03-12'
$ws.Cells.Item(39, 1).Value = 'WI'
$ws.Cells.Item(39, 2).Value = '$3.00 Games'
$ws.Cells.Item(39, 3).Value = 'Magic Word Crossword'
$ws.Cells.Item(39, 4).Value = 2089
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = '''2019-03-12'
$ws.Cells.Item(40, 6).Value = '''2019-03-12'
$ws.Cells.Item(41, 6).Value = '''2019-03-12'
$ws.Cells.Item(42, 1).Value = 'WI'
$ws.Cells.Item(42, 2).Value = '$3.00 Games'
$ws.Cells.Item(42, 3).Value = 'Polka Crossword'
$ws.Cells.Item(42, 4).Value = 2119
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = '''2019-03-12'
$ws.Cells.Item(43, 1).Value = 'WI'
$ws.Cells.Item(43, 2).Value = '$3.00 Games'
$ws.Cells.Item(43, 3).Value = 'Naughty or Nice Crossword'
$ws.Cells.Item(43, 4).Value = 2116
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = '''2019-03-12'
$ws.Cells.Item(44, 1).Value = 'WI'
$ws.Cells.Item(44, 2).Value = '$3.00 Games'
$ws.Cells.Item(44, 3).Value = 'Crossword Craving'
$ws.Cells.Item(44, 4).Value = 2081
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = '''2019-03-12'
$ws.Cells.Item(45, 1).Value = 'WI'
$ws.Cells.Item(45, 2).Value = '$3.00 Games'
$ws.Cells.Item(45, 3).Value = 'Fruit Explosion'
$ws.Cells.Item(45, 4).Value = 2071
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = '''2019-03-12'
$ws.Cells.Item(46, 1).Value = 'WI'
$ws.Cells.Item(46, 2).Value = '$3.00 Games'
$ws.Cells.Item(46, 3).Value = 'Lucky Charm Slingo'
$ws.Cells.Item(46, 4).Value = 2046
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = '''2019-03-12'
$ws.Cells.Item(47, 1).Value = 'WI'
$ws.Cells.Item(47, 2).Value = '$3.00 Games'
$ws.Cells.Item(47, 3).Value = 'Here Bingo!'
$ws.Cells.Item(47, 4).Value = 2080
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = '''2019-03-12'
$ws.Cells.Item(48, 1).Value = 'WI'
$ws.Cells.Item(48, 2).Value = '$3.00 Games'
$ws.Cells.Item(48, 3).Value = 'Wild Bingo'
$ws.Cells.Item(48, 4).Value = 2049
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = '''2019-03-12'
$ws.Cells.Item(49, 1).Value = 'WI'
$ws.Cells.Item(49, 2).Value = '$3.00 Games'
$ws.Cells.Item(49, 3).Value = 'Badger State Slingo'
$ws.Cells.Item(49, 4).Value = 2074
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = '''2019-03-12'
$ws.Cells.Item(50, 1).Value = 'WI'
$ws.Cells.Item(50, 2).Value = '$5.00 Games'
$ws.Cells.Item(50, 3).Value = 'Lucky Letter Crossword'
$ws.Cells.Item(50, 4).Value = 2026
$ws.Cells.Item(50, 5).Value = 1
$ws.Cells.Item(50, 6).Value = '''2019-03-12'
$ws.Cells.Item(51, 1).Value = 'WI'
$ws.Cells.Item(51, 2).Value = '$5.00 Games'
$ws.Cells.Item(51, 3).Value = 'Road To Riches'
$ws.Cells.Item(51, 4).Value = 2139
$ws.Cells.Item(51, 5).Value = 1
$ws.Cells.Item(51, 6).Value = '''2019-03-12'
$ws.Cells.Item(52, 1).Value = 'WI'
$ws.Cells.Item(52, 2).Value = '$5.00 Games'
$ws.Cells.Item(52, 3).Value = '20X The Money'
$ws.Cells.Item(52, 4).Value = 2150
$ws.Cells.Item(52, 5).Value = 2
$ws.Cells.Item(52, 6).Value = '''2019-03-12'
$ws.Cells.Item(53, 1).Value = 'WI'
$ws.Cells.Item(53, 2).Value = '$5.00 Games'
$ws.Cells.Item(53, 3).Value = 'Super Mega Bonus '
$ws.Cells.Item(53, 4).Value = 2093
$ws.Cells.Item(53, 5).Value = 1
$ws.Cells.Item(53, 6).Value = '''2019-03-12'
$ws.Cells.Item(54, 1).Value = 'WI'
$ws.Cells.Item(54, 2).Value = '$5.00 Games'
$ws.Cells.Item(54, 3).Value = 'Mystery Letter Crossword'
$ws.Cells.Item(54, 4).Value = 2134
$ws.Cells.Item(54, 5).Value = 1
$ws.Cells.Item(54, 6).Value = '''2019-03-12'
$ws.Cells.Item(55, 1).Value = 'WI'
$ws.Cells.Item(55, 2).Value = '$5.00 Games'
$ws.Cells.Item(55, 3).Value = 777
$ws.Cells.Item(55, 4).Value = 2061
$ws.Cells.Item(55, 5).Value = 1
$ws.Cells.Item(55, 6).Value = '''2019-03-12'
$ws.Cells.Item(56, 1).Value = 'WI'
$ws.Cells.Item(56, 2).Value = '$5.00 Games'
$ws.Cells.Item(56, 3).Value = 'Extend Your Cash'
$ws.Cells.Item(56, 4).Value = 2123
$ws.Cells.Item(56, 5).Value = 1
$ws.Cells.Item(56, 6).Value = '''2019-03-12'
$ws.Cells.Item(57, 1).Value = 'WI'
$ws.Cells.Item(57, 2).Value = '$5.00 Games'
$ws.Cells.Item(57, 3).Value = 'Triple Play'
$ws.Cells.Item(57, 4).Value = 2048
$ws.Cells.Item(57, 5).Value = 1
$ws.Cells.Item(57, 6).Value = '''2019-03-12'
$ws.Cells.Item(58, 1).Value = 'WI'
$ws.Cells.Item(58, 2).Value = '$5.00 Games'
$ws.Cells.Item(58, 3).Value = 'Cashilicious'
$ws.Cells.Item(58, 4).Value = 2137
$ws.Cells.Item(58, 5).Value = 3
$ws.Cells.Item(58, 6).Value = '''2019-03-12'
$ws.Cells.Item(59, 1).Value = 'WI'
$ws.Cells.Item(59, 2).Value = '$5.00 Games'
$ws.Cells.Item(59, 3).Value = 'Supper Club Cash'
$ws.Cells.Item(59, 4).Value = 2125
$ws.Cells.Item(59, 5).Value = 469
$ws.Cells.Item(59, 6).Value = '''2019-03-12'
$ws.Cells.Item(60, 1).Value = 'WI'
$ws.Cells.Item(60, 2).Value = '$5.00 Games'
$ws.Cells.Item(60, 3).Value = 'Jolly $500''s '
$ws.Cells.Item(60, 4).Value = 2115
$ws.Cells.Item(60, 5).Value = 49
$ws.Cells.Item(60, 6).Value = '''2019-03-12'
$ws.Cells.Item(61, 1).Value = 'WI'
$ws.Cells.Item(61, 2).Value = '$5.00 Games'
$ws.Cells.Item(61, 3).Value = 'Land Of Milk And Money'
$ws.Cells.Item(61, 4).Value = 2090
$ws.Cells.Item(61, 5).Value = 27
$ws.Cells.Item(61, 6).Value = '''2019-03-12'
$ws.Cells.Item(62, 1).Value = 'WI'
$ws.Cells.Item(62, 2).Value = '$5.00 Games'
$ws.Cells.Item(62, 3).Value = 'In The Green'
$ws.Cells.Item(62, 4).Value = 2091
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(62, 6).Value = '''2019-03-12'
$ws.Cells.Item(63, 1).Value = 'WI'
$ws.Cells.Item(63, 2).Value = '$5.00 Games'
$ws.Cells.Item(63, 3).Value = 'Platinum Crossword'
$ws.Cells.Item(63, 4).Value = 2054
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = '''2019-03-12'
$ws.Cells.Item(64, 1).Value = 'WI'
$ws.Cells.Item(64, 2).Value = '$5.00 Games'
$ws.Cells.Item(64, 3).Value = '5 Star Crossword'
$ws.Cells.Item(64, 4).Value = 2109
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = '''2019-03-12'
$ws.Cells.Item(65, 1).Value = 'WI'
$ws.Cells.Item(65, 2).Value = '$5.00 Games'
$ws.Cells.Item(65, 3).Value = 'Pack Attack'
$ws.Cells.Item(65, 4).Value = 2104
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = '''2019-03-12'
$ws.Cells.Item(66, 1).Value = 'WI'
$ws.Cells.Item(66, 2).Value = '$5.00 Games'
$ws.Cells.Item(66, 3).Value = 'Wild Cherry Crossword'
$ws.Cells.Item(66, 4).Value = 2082
$ws.Cells.Item(66, 5).Value = 0
$ws.Cells.Item(66, 6).Value = '''2019-03-12'
$ws.Cells.Item(67, 1).Value = 'WI'
$ws.Cells.Item(67, 2).Value = '$5.00 Games'
$ws.Cells.Item(67, 3).Value = '30 Years of Fun!'
$ws.Cells.Item(67, 4).Value = 2058
$ws.Cells.Item(67, 5).Value = 0
$ws.Cells.Item(67, 6).Value = '''2019-03-12'
$ws.Cells.Item(68, 1).Value = 'WI'
$ws.Cells.Item(68, 2).Value = '$5.00 Games'
$ws.Cells.Item(68, 3).Value = 'Cash Up'
$ws.Cells.Item(68, 4).Value = 2064
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = '''2019-03-12'
$ws.Cells.Item(69, 1).Value = 'WI'
$ws.Cells.Item(69, 2).Value = '$5.00 Games'
$ws.Cells.Item(69, 3).Value = 'Vegas Nights'
$ws.Cells.Item(69, 4).Value = 2070
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = '''2019-03-12'
$ws.Cells.Item(70, 1).Value = 'WI'
$ws.Cells.Item(70, 2).Value = '$5.00 Games'
$ws.Cells.Item(70, 3).Value = 'Boost Your Bucks'
$ws.Cells.Item(70, 4).Value = 2122
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = '''2019-03-12'
$ws.Cells.Item(71, 1).Value = 'WI'
$ws.Cells.Item(71, 2).Value = '$5.00 Games'
$ws.Cells.Item(71, 3).Value = 'Gold Rush'
$ws.Cells.Item(71, 4).Value = 2023
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = '''2019-03-12'
$ws.Cells.Item(72, 1).Value = 'WI'
$ws.Cells.Item(72, 2).Value = '$5.00 Games'
$ws.Cells.Item(72, 3).Value = 'Deluxe Bucks'
$ws.Cells.Item(72, 4).Value = 2092
$ws.Cells.Item(72, 5).Value = 0
$ws.Cells.Item(72, 6).Value = '''2019-03-12'
$ws.Cells.Item(73, 1).Value = 'WI'
$ws.Cells.Item(73, 2).Value = '$10.00 Games'
$ws.Cells.Item(73, 3).Value = 'Double Dollar Crossword'
$ws.Cells.Item(73, 4).Value = 2154
$ws.Cells.Item(73, 5).Value = 3
$ws.Cells.Item(73, 6).Value = '''2019-03-12'
$ws.Cells.Item(74, 1).Value = 'WI'
$ws.Cells.Item(74, 2).Value = '$10.00 Games'
$ws.Cells.Item(74, 3).Value = 'Winner''s Circle'
$ws.Cells.Item(74, 4).Value = 2094
$ws.Cells.Item(74, 5).Value = 2
$ws.Cells.Item(74, 6).Value = '''2019-03-12'
$ws.Cells.Item(75, 1).Value = 'WI'
$ws.Cells.Item(75, 2).Value = '$10.00 Games'
$ws.Cells.Item(75, 3).Value = '50X The Money'
$ws.Cells.Item(75, 4).Value = 2151
$ws.Cells.Item(75, 5).Value = 2
$ws.Cells.Item(75, 6).Value = '''2019-03-12'
$ws.Cells.Item(76, 1).Value = 'WI'
$ws.Cells.Item(76, 2).Value = '$10.00 Games'
$ws.Cells.Item(76, 3).Value = 'Hot $500''s'
$ws.Cells.Item(76, 4).Value = 2124
$ws.Cells.Item(76, 5).Value = 22
$ws.Cells.Item(76, 6).Value = '''2019-03-12'
$ws.Cells.Item(77, 1).Value = 'WI'
$ws.Cells.Item(77, 2).Value = '$10.00 Games'
$ws.Cells.Item(77, 3).Value = 'Finding $500''s'
$ws.Cells.Item(77, 4).Value = 2055
$ws.Cells.Item(77, 5).Value = 19
$ws.Cells.Item(77, 6).Value = '''2019-03-12'
$ws.Cells.Item(78, 1).Value = 'WI'
$ws.Cells.Item(78, 2).Value = '$10.00 Games'
$ws.Cells.Item(78, 3).Value = 'Full of $500''s'
$ws.Cells.Item(78, 4).Value = 2157
$ws.Cells.Item(78, 5).Value = 2641
$ws.Cells.Item(78, 6).Value = '''2019-03-12'
$ws.Cells.Item(79, 1).Value = 'WI'
$ws.Cells.Item(79, 2).Value = '$10.00 Games'
$ws.Cells.Item(79, 3).Value = 'EZ Grand'
$ws.Cells.Item(79, 4).Value = 2135
$ws.Cells.Item(79, 5).Value = 154
$ws.Cells.Item(79, 6).Value = '''2019-03-12'
$ws.Cells.Item(80, 1).Value = 'WI'
$ws.Cells.Item(80, 2).Value = '$10.00 Games'
$ws.Cells.Item(80, 3).Value = 'Hit $1,000'
$ws.Cells.Item(80, 4).Value = 2043
$ws.Cells.Item(80, 5).Value = 15
$ws.Cells.Item(80, 6).Value = '''2019-03-12'
$ws.Cells.Item(81, 1).Value = 'WI'
$ws.Cells.Item(81, 2).Value = '$10.00 Games'
$ws.Cells.Item(81, 3).Value = 'Instant $1,000'
$ws.Cells.Item(81, 4).Value = 2087
$ws.Cells.Item(81, 5).Value = 19
$ws.Cells.Item(81, 6).Value = '''2019-03-12'
$ws.Cells.Item(82, 1).Value = 'WI'
$ws.Cells.Item(82, 2).Value = '$10.00 Games'
$ws.Cells.Item(82, 3).Value = 'All $50''s and $100''s'
$ws.Cells.Item(82, 4).Value = 2121
$ws.Cells.Item(82, 5).Value = 18446
$ws.Cells.Item(82, 6).Value = '''2019-03-12'
$ws.Cells.Item(83, 1).Value = 'WI'
$ws.Cells.Item(83, 2).Value = '$10.00 Games'
$ws.Cells.Item(83, 3).Value = '$50''s and $100''s'
$ws.Cells.Item(83, 4).Value = 2066
$ws.Cells.Item(83, 5).Value = 170
$ws.Cells.Item(83, 6).Value = '''2019-03-12'
$ws.Cells.Item(84, 1).Value = 'WI'
$ws.Cells.Item(84, 2).Value = '$10.00 Games'
$ws.Cells.Item(84, 3).Value = 'Incredible Crossword'
$ws.Cells.Item(84, 4).Value = 2052
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = '''2019-03-12'
$ws.Cells.Item(85, 1).Value = 'WI'
$ws.Cells.Item(85, 2).Value = '$10.00 Games'
$ws.Cells.Item(85, 3).Value = 'Champion Crossword'
$ws.Cells.Item(85, 4).Value = 2126
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = '''2019-03-12'
$ws.Cells.Item(86, 1).Value = 'WI'
$ws.Cells.Item(86, 2).Value = '$10.00 Games'
$ws.Cells.Item(86, 3).Value = 'Cash Craze Crossword '
$ws.Cells.Item(86, 4).Value = 2107
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = '''2019-03-12'
$ws.Cells.Item(87, 1).Value = 'WI'
$ws.Cells.Item(87, 2).Value = '$10.00 Games'
$ws.Cells.Item(87, 3).Value = 'Joy'
$ws.Cells.Item(87, 4).Value = 2114
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = '''2019-03-12'
$ws.Cells.Item(88, 1).Value = 'WI'
$ws.Cells.Item(88, 2).Value = '$10.00 Games'
$ws.Cells.Item(88, 3).Value = 'Perfect 10'
$ws.Cells.Item(88, 4).Value = 2084
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = '''2019-03-12'
$ws.Cells.Item(89, 1).Value = 'WI'
$ws.Cells.Item(89, 2).Value = '$10.00 Games'
$ws.Cells.Item(89, 3).Value = 'Stinkin'' Rich'
$ws.Cells.Item(89, 4).Value = 2086
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = '''2019-03-12'
$ws.Cells.Item(90, 1).Value = 'WI'
$ws.Cells.Item(90, 2).Value = '$10.00 Games'
$ws.Cells.Item(90, 3).Value = '100 Seasons'
$ws.Cells.Item(90, 4).Value = 2105
$ws.Cells.Item(90, 5).Value = 0
$ws.Cells.Item(90, 6).Value = '''2019-03-12'
$ws.Cells.Item(91, 1).Value = 'WI'
$ws.Cells.Item(91, 2).Value = '$15.00 Games'
$ws.Cells.Item(91, 3).Value = 'Holiday Countdown '
$ws.Cells.Item(91, 4).Value = 2113
$ws.Cells.Item(91, 5).Value = 0
$ws.Cells.Item(91, 6).Value = '''2019-03-12'
$ws.Cells.Item(92, 1).Value = 'WI'
$ws.Cells.Item(92, 2).Value = '$20.00 Games'
$ws.Cells.Item(92, 3).Value = 'Hit It Big'
$ws.Cells.Item(92, 4).Value = 2112
$ws.Cells.Item(92, 5).Value = 1
$ws.Cells.Item(92, 6).Value = '''2019-03-12'
$ws.Cells.Item(93, 1).Value = 'WI'
$ws.Cells.Item(93, 2).Value = '$20.00 Games'
$ws.Cells.Item(93, 3).Value = 'Power Up Your Crossword'
$ws.Cells.Item(93, 4).Value = 2131
$ws.Cells.Item(93, 5).Value = 1
$ws.Cells.Item(93, 6).Value = '''2019-03-12'
$ws.Cells.Item(94, 1).Value = 'WI'
$ws.Cells.Item(94, 2).Value = '$20.00 Games'
$ws.Cells.Item(94, 3).Value = 'Jackpot Payout '
$ws.Cells.Item(94, 4).Value = 2111
$ws.Cells.Item(94, 5).Value = 1
$ws.Cells.Item(94, 6).Value = '''2019-03-12'
$ws.Cells.Item(95, 1).Value = 'WI'
$ws.Cells.Item(95, 2).Value = '$20.00 Games'
$ws.Cells.Item(95, 3).Value = '$100,000 Large'
$ws.Cells.Item(95, 4).Value = 2133
$ws.Cells.Item(95, 5).Value = 3
$ws.Cells.Item(95, 6).Value = '''2019-03-12'
$ws.Cells.Item(96, 1).Value = 'WI'
$ws.Cells.Item(96, 2).Value = '$20.00 Games'
$ws.Cells.Item(96, 3).Value = 'Cafe Crossword'
$ws.Cells.Item(96, 4).Value = 2088
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(96, 6).Value = '''2019-03-12'
$ws.Cells.Item(97, 1).Value = 'WI'
$ws.Cells.Item(97, 2).Value = '$20.00 Games'
$ws.Cells.Item(97, 3).Value = '$200,000 Cash Reserve'
$ws.Cells.Item(97, 4).Value = 2075
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = '''2019-03-12'
$ws.Cells.Item(98, 1).Value = 'WI'
$ws.Cells.Item(98, 2).Value = '$20.00 Games'
$ws.Cells.Item(98, 3).Value = 'Mega Crossword'
$ws.Cells.Item(98, 4).Value = 2067
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = '''2019-03-12'
$ws.Cells.Item(99, 1).Value = 'WI'
$ws.Cells.Item(99, 2).Value = '$20.00 Games'
$ws.Cells.Item(99, 3).Value = 'Pot Of Gold Crossword'
$ws.Cells.Item(99, 4).Value = 2118
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = '''2019-03-12'
$ws.Cells.Item(100, 1).Value = 'WI'
$ws.Cells.Item(100, 2).Value = '$20.00 Games'
$ws.Cells.Item(100, 3).Value = '$200 Grand'
$ws.Cells.Item(100, 4).Value = 2044
$ws.Cells.Item(100, 5).Value = 0
$ws.Cells.Item(100, 6).Value = '''2019-03-12'
$ws.Cells.Item(101, 1).Value = 'WI'
$ws.Cells.Item(101, 2).Value = '$20.00 Games'
$ws.Cells.Item(101, 3).Value = 'Epic Multiplier'
$ws.Cells.Item(101, 4).Value = 2083
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = '''2019-03-12'
$ws.Cells.Item(102, 1).Value = 'WI'
$ws.Cells.Item(102, 2).Value = '$20.00 Games'
$ws.Cells.Item(102, 3).Value = 'Cash Attack'
$ws.Cells.Item(102, 4).Value = 2047
$ws.Cells.Item(102, 5).Value = 0
$ws.Cells.Item(102, 6).Value = '''2019-03-12'
$ws.Cells.Item(103, 1).Value = 'WI'
$ws.Cells.Item(103, 2).Value = '$30.00 Games'
$ws.Cells.Item(103, 3).Value = '$10,000 Cash'
$ws.Cells.Item(103, 4).Value = 655
$ws.Cells.Item(103, 5).Value = 26
$ws.Cells.Item(103, 6).Value = '''2019-03-12'
$ws.Cells.Item(104, 1).Value = 'WI'
$ws.Cells.Item(104, 2).Value = '$30.00 Games'
$ws.Cells.Item(104, 3).Value = 'Golden Millions'
$ws.Cells.Item(104, 4).Value = 2130
$ws.Cells.Item(104, 5).Value = 2
$ws.Cells.Item(104, 6).Value = '''2019-03-12'
$ws.Cells.Item(105, 1).Value = 'WI'
$ws.Cells.Item(105, 2).Value = '$30.00 Games'
$ws.Cells.Item(105, 3).Value = 'Magnificent Millions'
$ws.Cells.Item(105, 4).Value = 2106
$ws.Cells.Item(105, 5).Value = 1
$ws.Cells.Item(105, 6).Value = '''2019-03-12'
$ws.Cells.Item(106, 1).Value = 'WI'
$ws.Cells.Item(106, 2).Value = '$30.00 Games'
$ws.Cells.Item(106, 3).Value = 'Wisconsin Millions'
$ws.Cells.Item(106, 4).Value = 2021
$ws.Cells.Item(106, 5).Value = 0
$ws.Cells.Item(106, 6).Value = '''2019-03-12'
$ws.Cells.Item(107, 1).Value = 'WI'
$ws.Cells.Item(107, 2).Value = '$30.00 Games'
$ws.Cells.Item(107, 3).Value = 'One Million Now!'
$ws.Cells.Item(107, 4).Value = 2057
$ws.Cells.Item(107, 5).Value = 0
$ws.Cells.Item(107, 6).Value = '''2019-03-12'
